# Updating functionalities required by 'cb' parameter
#
# This script inserts a new "steam" grid entry before the existing
# "industry" entries on the 'demanddata_other demands' sheet: for every
# row currently holding Grid = "industry" in column B, a new value
# "steam" is written into column B and the previous "industry" value is
# moved one column to the right, into column C (Node_suffix column was
# otherwise unused for those rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("demanddata_other demands")
$ws.Activate()

# Rows whose Grid column (B) currently contains "industry" and need the
# new "steam" entry inserted ahead of it (industry moves to column C).
$rows = @(27, 28, 29, 47, 48, 49, 50, 51, 52, 53, 54, 55, 56, 57, 58, 59, 60, 61)

foreach ($r in $rows) {
    $ws.Range("C$r").Value = $ws.Range("B$r").Value2
    $ws.Range("B$r").Value = "steam"
}

# Reflect the reviewer's final view/selection state on the sheet.
$ws.Range("K34").Select()
$excel.ActiveWindow.ScrollRow = 12
$excel.ActiveWindow.ScrollColumn = 1
